# Applies the "Added GRU results, slight improvements" edit:
#  - row 5: G5/H5/I5 get new text values (101/102/103 embedding dim)
#  - new row 6 added with GRU results
#  - (font for phonetic info is created implicitly by Excel when it needs one;
#    no direct COM surface for that, so we focus on the data change)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- Add new row 6 with GRU results (order matters for shared-string table) ---
$ws.Range("E6").Value = "128 GRU units"

# --- Update existing row 5 cells (text replaced by new shared strings) ---
$ws.Range("G5").Value = "101 embedding dim"
$ws.Range("H5").Value = "102 embedding dim"
$ws.Range("I5").Value = "103 embedding dim"

# --- Finish filling new row 6 ---
$ws.Range("F6").Value = "100 embedding dim"
$ws.Range("G6").Value = "101 embedding dim"
$ws.Range("H6").Value = "102 embedding dim"
$ws.Range("I6").Value = "103 embedding dim"
$ws.Range("J6").Value = "~"
$ws.Range("K6").Value = "{'exact_match': 12.039666651147632, 'f1': 14.165593903763927}"

# Match formatting of row 5 (style index 1 => wrapText, row height 45)
$ws.Range("E6:K6").Style = $ws.Range("E5:K5").Style
$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(5).RowHeight

# Update selection to mirror the saved workbook state (J6 selected)
$ws.Range("J6").Select()
